$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($ws, $row1, $row2, $startCol, $endCol) {
    for ($c = $startCol; $c -le $endCol; $c++) {
        $v1 = $ws.Cells.Item($row1, $c).Value()
        $v2 = $ws.Cells.Item($row2, $c).Value()
        $ws.Cells.Item($row1, $c).Value = $v2
        $ws.Cells.Item($row2, $c).Value = $v1
    }
}

# Swap match data (columns F:V) between row pairs whose order was corrected
Swap-RowRange $ws 24 25 6 22
Swap-RowRange $ws 26 29 6 22
Swap-RowRange $ws 43 44 6 22

# Add new row 46: Palermo vs Cosenza
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "italy"
$ws.Cells.Item(46, 3).Value = "serie-b"
$ws.Cells.Item(46, 4).Value = "2023-2024"
$ws.Cells.Item(46, 5).Value = 45191.85416666666
$ws.Cells.Item(46, 6).Value = "Palermo"
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = "Cosenza"
$ws.Cells.Item(46, 9).Value = 1
$ws.Cells.Item(46, 10).Value = 1.68
$ws.Cells.Item(46, 11).Value = "16/09/2023 13:14"
$ws.Cells.Item(46, 12).Value = 1.74
$ws.Cells.Item(46, 13).Value = "22/09/2023 20:22"
$ws.Cells.Item(46, 14).Value = 3.94
$ws.Cells.Item(46, 15).Value = "16/09/2023 13:14"
$ws.Cells.Item(46, 16).Value = 3.68
$ws.Cells.Item(46, 17).Value = "22/09/2023 20:28"
$ws.Cells.Item(46, 18).Value = 5.2
$ws.Cells.Item(46, 19).Value = "16/09/2023 13:14"
$ws.Cells.Item(46, 20).Value = 5.5
$ws.Cells.Item(46, 21).Value = "22/09/2023 20:28"
$ws.Cells.Item(46, 22).Value = "https://www.betexplorer.com/football/italy/serie-b/palermo-cosenza/jRjNmct8/"

# Match the formatting used by the other data rows (index style + date style)
$ws.Cells.Item(45, 1).Copy() | Out-Null
$ws.Cells.Item(46, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(45, 5).Copy() | Out-Null
$ws.Cells.Item(46, 5).PasteSpecial(-4122) | Out-Null

# restore the value after the format paste (PasteSpecial formats only, but be safe)
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 5).Value = 45191.85416666666
